$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.012.81'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '2.461.81'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'519.88"
$ws.Range("E5").Value = '  -2.66%  '
$ws.Range("D6").Value = "'133.17"
$ws.Range("E6").Value = '  -2.46%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.43%  '
$ws.Range("D8").Value = "'0.556"
$ws.Range("E8").Value = '  -1.45%  '
$ws.Range("D9").Value = '2.469.34'
$ws.Range("E9").Value = '  -1.82%  '
$ws.Range("D10").Value = "'0.0980"
$ws.Range("E10").Value = '  -3.20%  '
$ws.Range("E11").Value = '  -0.53%  '
$ws.Range("D12").Value = "'5.26"
$ws.Range("E12").Value = '  -2.90%  '
$ws.Range("D13").Value = "'0.339"
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("D14").Value = '2.898.57'
$ws.Range("E14").Value = '  -2.18%  '
$ws.Range("D15").Value = '57.945.95'
$ws.Range("E15").Value = '  -1.49%  '
$ws.Range("D16").Value = "'22.29"
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").Value = "'0.0000135"
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("D18").Value = '2.470.90'
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").Value = "'10.62"
$ws.Range("E19").Value = '  -4.16%  '
$ws.Range("D20").Value = "'320.19"
$ws.Range("E20").Value = '  -0.93%  '
$ws.Range("D21").Value = "'4.15"
$ws.Range("E21").Value = '  -2.41%  '
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = '  -0.10%  '
$ws.Range("D23").Value = "'5.72"
$ws.Range("E23").Value = '  -4.41%  '
$ws.Range("D24").Value = "'64.57"
$ws.Range("E24").Value = '  -0.94%  '
$ws.Range("D25").Value = "'0.408"
$ws.Range("E25").Value = '  -2.76%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("D27").Value = "'0.159"
$ws.Range("E27").Value = '  -3.02%  '
$ws.Range("D28").Value = "'7.34"
$ws.Range("E28").Value = '  -2.36%  '
$ws.Range("D29").Value = '0.0₃0746'
$ws.Range("E29").Value = '  -2.86%  '
$ws.Range("D30").Value = "'167.85"
$ws.Range("E30").Value = '  -2.25%  '
$ws.Range("D31").Value = "'1.69"
$ws.Range("E31").Value = '  -3.42%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").Value = "'6.22"
$ws.Range("E32").Value = '  -5.79%  '
$ws.Range("B33").Value = 'Fetch.AI'
$ws.Range("C33").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D33").Value = "'1.17"
$ws.Range("E33").Value = '  -0.48%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").Value = "'18.04"
$ws.Range("E37").Value = '  -1.76%  '
$ws.Range("D38").Value = "'3.96"
$ws.Range("E38").Value = '  -2.38%  '
$ws.Range("D39").Value = "'1.47"
$ws.Range("E39").Value = '  -4.51%  '
$ws.Range("D40").Value = "'36.29"
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").Value = "'0.798"
$ws.Range("E41").Value = '  -1.21%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = "'274.82"
$ws.Range("E42").Value = '  -3.27%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").Value = "'3.45"
$ws.Range("E43").Value = '  -3.51%  '
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").Value = "'5.03"
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").Value = "'0.589"
$ws.Range("E45").Value = '  -3.07%  '
$ws.Range("D46").Value = "'124.00"
$ws.Range("E46").Value = '  -4.63%  '
$ws.Range("D47").Value = "'0.0909"
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").Value = "'0.0489"
$ws.Range("E48").Value = '  -2.96%  '
$ws.Range("D49").Value = "'0.0212"
$ws.Range("E49").Value = '  -2.87%  '
$ws.Range("D50").Value = "'16.86"
$ws.Range("E50").Value = '  -2.76%  '
$ws.Range("D51").Value = '1.724.17'
$ws.Range("E51").Value = '  -1.85%  '
